# Daily attendance processing - 2025-11-22 16:24:50
# Normalize the "Recorded By" (column G) entries so that the ordering of
# names/emails listed in each cell follows the corrected convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match replacements applied to the "Recorded By" column values.
$replacements = @{
    "backup@backdoor.com, system, System" = "system, backup@backdoor.com, System"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
